$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (the old "No. Nómina" / "Nombre(s)" headers were
# replaced with shorter "NumNomina" / "Nombres" labels).
$ws.Range("A1").Value = "NumNomina"
$ws.Range("B1").Value = "Nombres"

# Shorten the second teacher's first name.
$ws.Range("B2").Value = "Elizatbeth"

# Move the active selection to G6, as left by the author when saving.
[void]$ws.Range("G6").Select()
